$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (date) values between rows 2-3 and rows 6-7.
$ws.Range("D2").Value = 44838
$ws.Range("D3").Value = 44838
$ws.Range("D6").Value = 44832
$ws.Range("D7").Value = 44832
